# Fruta / hortaliza, semanal
#
# The weekly refresh prepends two new daily-price records (dated 2022-01-31,
# serial 44592) to the Mango subset sheet. Every existing record (rows
# 36-124) shifts down two rows (to rows 38-126); the dimension grows from
# A1:T124 to A1:T126.
#
# Implementation: insert two blank rows at row 36 (only within columns A:T,
# so we don't balloon the row width to the full 16384-column sheet), then
# seed those two new rows with a copy of the data that is now directly
# below them (the shifted former row 36 / row 37), and finally overwrite
# the date (column D) of the two new rows with the new date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 36, shifting rows 36:124 down to 38:126.
$ws.Range("A36:T37").Insert()

# The two new (currently blank) rows should start life as copies of the
# rows now sitting right underneath them (the shifted originals), then get
# the new date stamped over column D.
$ws.Range("A38:T38").Copy($ws.Range("A36:T36"))
$ws.Range("A39:T39").Copy($ws.Range("A37:T37"))

$ws.Range("D36").Value = 44592
$ws.Range("D37").Value = 44592
